$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 7 to make room for the IBF Financial Holdings row,
# shifting the existing "Hotai Finance" row (old row 7) down to row 8.
$ws.Rows.Item(7).Insert()

# --- Row 2 ---
$ws.Cells.Item(2, 1).Value = "Taiwan"
# B2 holds a numeric-looking company identifier that must stay a text string
# (matches the source data's inline-string type), so force text formatting
# before assigning, then restore the cell's (unstyled) appearance.
$ws.Cells.Item(2, 2).NumberFormat = "@"
$ws.Cells.Item(2, 2).Value = "6"
$ws.Cells.Item(2, 2).Style = "Normal"
$ws.Cells.Item(2, 3).Value = "Financial Svcs. (Non-bank & Insurance)"
$ws.Cells.Item(2, 4).Value = 0.09210000000000002
$ws.Cells.Item(2, 5).Value = 0.151
$ws.Cells.Item(2, 7).Value = 0.05535205535205535
$ws.Cells.Item(2, 8).Value = 0.05423722815027163
$ws.Cells.Item(2, 9).Value = 0.02213551343986126
$ws.Cells.Item(2, 10).Value = 0.01746022222379003
$ws.Cells.Item(2, 11).Value = 1646.59
$ws.Cells.Item(2, 12).Value = 0.2913751305055652
$ws.Cells.Item(2, 13).Value = 681.7889
$ws.Cells.Item(2, 14).Value = 0.03088959214925833
$ws.Cells.Item(2, 15).Value = 0.4140611202545868
$ws.Cells.Item(2, 16).Value = 664.4189
$ws.Cells.Item(2, 17).Value = 0.03010261510162289
$ws.Cells.Item(2, 18).Value = 0.4035120461074099
$ws.Cells.Item(2, 19).Value = 17.36999999999999
$ws.Cells.Item(2, 20).Value = 0.02547709415627035
$ws.Cells.Item(2, 21).Value = 4156.4
$ws.Cells.Item(2, 22).Value = 0.1883126885890593
$ws.Cells.Item(2, 23).Value = 0.2121902107987375
$ws.Cells.Item(2, 24).Value = 0.03411710226802814
$ws.Cells.Item(2, 25).Value = 0.1780731085307093
$ws.Cells.Item(2, 26).Value = 0.1153282442729397
$ws.Cells.Item(2, 27).Value = 0
$ws.Cells.Item(2, 28).Value = 0.02430309416525871
$ws.Cells.Item(2, 29).Value = -0.02332028293332264
$ws.Cells.Item(2, 30).Value = 47272.305
$ws.Cells.Item(2, 31).Value = 0
$ws.Cells.Item(2, 32).Value = 47272.305
$ws.Cells.Item(2, 33).Value = 43115.905
$ws.Cells.Item(2, 34).Value = 0.6817061810805692
$ws.Cells.Item(2, 35).Value = 0.7547157656618836
$ws.Cells.Item(2, 36).Value = 0.661411611284674
$ws.Cells.Item(2, 37).Value = 0.7372823179676367
$ws.Cells.Item(2, 38).Value = 0.166
$ws.Cells.Item(2, 39).Value = -8.238999999999999
$ws.Cells.Item(2, 40).Value = 151.4458416095342
$ws.Cells.Item(2, 41).Value = 753.55421686747
$ws.Cells.Item(2, 42).Value = 138.1300217850964
$ws.Cells.Item(2, 43).Value = -15.18266779949023

# --- Row 3 ---
$ws.Cells.Item(3, 1).Value = "Taiwan"
$ws.Cells.Item(3, 2).Value = "SysJust Co., Ltd. (GTSM:3158)"
$ws.Cells.Item(3, 3).Value = "Financial Svcs. (Non-bank & Insurance)"
$ws.Cells.Item(3, 4).Value = 0.0495
$ws.Cells.Item(3, 5).Value = 0.0127
$ws.Cells.Item(3, 7).Value = 0.4523809523809524
$ws.Cells.Item(3, 8).Value = 0.1523809523809524
$ws.Cells.Item(3, 9).Value = 0.1566666666666667
$ws.Cells.Item(3, 10).Value = 0.1304776119402985
$ws.Cells.Item(3, 11).Value = 2.79
$ws.Cells.Item(3, 12).Value = 0.1328571428571429
$ws.Cells.Item(3, 13).Value = 2.23
$ws.Cells.Item(3, 14).Value = 0.0501123595505618
$ws.Cells.Item(3, 15).Value = 0.7992831541218638
$ws.Cells.Item(3, 16).Value = 2.23
$ws.Cells.Item(3, 17).Value = 0.0501123595505618
$ws.Cells.Item(3, 18).Value = 0.7992831541218638
$ws.Cells.Item(3, 19).Value = 0
$ws.Cells.Item(3, 20).Value = 0
$ws.Cells.Item(3, 21).Value = 5.5
$ws.Cells.Item(3, 22).Value = 0.1235955056179775
$ws.Cells.Item(3, 23).Value = 0.2196850393700787
$ws.Cells.Item(3, 24).Value = 0.0180597275576015
$ws.Cells.Item(3, 25).Value = 0.2016253118124773
$ws.Cells.Item(3, 26).Value = 2.349256068911511
$ws.Cells.Item(3, 27).Value = 0.3065253217078274
$ws.Cells.Item(3, 28).Value = 0.01807043912548618
$ws.Cells.Item(3, 29).Value = 0.2884548825823412
$ws.Cells.Item(3, 30).Value = 0.305
$ws.Cells.Item(3, 31).Value = 0
$ws.Cells.Item(3, 32).Value = 0.305
$ws.Cells.Item(3, 33).Value = -5.195
$ws.Cells.Item(3, 34).Value = 0.006807275973663654
$ws.Cells.Item(3, 35).Value = 0.02177793645126741
$ws.Cells.Item(3, 36).Value = -0.13217147945554
$ws.Cells.Item(3, 37).Value = -0.610817166372722
$ws.Cells.Item(3, 38).Value = 0.02
$ws.Cells.Item(3, 39).Value = -0.045
$ws.Cells.Item(3, 40).Value = 0.07549504950495049
$ws.Cells.Item(3, 41).Value = 164.5
$ws.Cells.Item(3, 42).Value = -1.285891089108911
$ws.Cells.Item(3, 43).Value = -73.11111111111111

# --- Row 4 ---
$ws.Cells.Item(4, 1).Value = "Taiwan"
$ws.Cells.Item(4, 2).Value = "Yulon Finance Corporation (TSEC:9941)"
$ws.Cells.Item(4, 3).Value = "Financial Svcs. (Non-bank & Insurance)"
$ws.Cells.Item(4, 4).Value = 0.116
$ws.Cells.Item(4, 5).Value = 0.155
$ws.Cells.Item(4, 7).Value = 0.305161485058859
$ws.Cells.Item(4, 8).Value = 0.305161485058859
$ws.Cells.Item(4, 9).Value = 0.1225475399939632
$ws.Cells.Item(4, 10).Value = 0.09609163313327722
$ws.Cells.Item(4, 11).Value = 105.5
$ws.Cells.Item(4, 12).Value = 0.1061474997484656
$ws.Cells.Item(4, 13).Value = 60.59999999999999
$ws.Cells.Item(4, 14).Value = 0.0509757738896366
$ws.Cells.Item(4, 15).Value = 0.5744075829383886
$ws.Cells.Item(4, 16).Value = 46.9
$ws.Cells.Item(4, 17).Value = 0.03945154777927322
$ws.Cells.Item(4, 18).Value = 0.4445497630331753
$ws.Cells.Item(4, 19).Value = 13.7
$ws.Cells.Item(4, 20).Value = 0.226072607260726
$ws.Cells.Item(4, 21).Value = 252.9
$ws.Cells.Item(4, 22).Value = 0.212735531628533
$ws.Cells.Item(4, 23).Value = 0.2046953822273962
$ws.Cells.Item(4, 24).Value = 0.05067887208046057
$ws.Cells.Item(4, 25).Value = 0.1540165101469356
$ws.Cells.Item(4, 26).Value = 0.1822933861560471
$ws.Cells.Item(4, 27).Value = 0.01751686918512971
$ws.Cells.Item(4, 28).Value = 0.02509224475030496
$ws.Cells.Item(4, 29).Value = -0.007575375565175253
$ws.Cells.Item(4, 30).Value = 5572
$ws.Cells.Item(4, 31).Value = 0
$ws.Cells.Item(4, 32).Value = 5572
$ws.Cells.Item(4, 33).Value = 5319.1
$ws.Cells.Item(4, 34).Value = 0.8241628209679328
$ws.Cells.Item(4, 35).Value = 0.8834908352888945
$ws.Cells.Item(4, 36).Value = 0.8173297069715268
$ws.Cells.Item(4, 37).Value = 0.8786236971208642
$ws.Cells.Item(4, 38).Value = 0
$ws.Cells.Item(4, 39).Value = 0
$ws.Cells.Item(4, 40).Value = 18.08503732554365
$ws.Cells.Item(4, 42).Value = 17.26419993508601

# --- Row 5 ---
$ws.Cells.Item(5, 1).Value = "Taiwan"
$ws.Cells.Item(5, 2).Value = "Chailease Holding Company Limited (TSEC:5871)"
$ws.Cells.Item(5, 3).Value = "Financial Svcs. (Non-bank & Insurance)"
$ws.Cells.Item(5, 4).Value = 0.0779
$ws.Cells.Item(5, 5).Value = 0.19
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 565.1
$ws.Cells.Item(5, 12).Value = 0.3653348849237135
$ws.Cells.Item(5, 13).Value = 234.2889
$ws.Cells.Item(5, 14).Value = 0.02558437346437346
$ws.Cells.Item(5, 15).Value = 0.4145972394266501
$ws.Cells.Item(5, 16).Value = 234.2889
$ws.Cells.Item(5, 17).Value = 0.02558437346437346
$ws.Cells.Item(5, 18).Value = 0.4145972394266501
$ws.Cells.Item(5, 19).Value = 0
$ws.Cells.Item(5, 20).Value = 0
$ws.Cells.Item(5, 21).Value = 1126.2
$ws.Cells.Item(5, 22).Value = 0.122981162981163
$ws.Cells.Item(5, 23).Value = 0.2644980107652703
$ws.Cells.Item(5, 24).Value = 0.0282664880665113
$ws.Cells.Item(5, 25).Value = 0.236231522698759
$ws.Cells.Item(5, 26).Value = 0.1268916228395077
$ws.Cells.Item(5, 27).Value = 0
$ws.Cells.Item(5, 28).Value = 0.02312662228643282
$ws.Cells.Item(5, 29).Value = -0.02312662228643282
$ws.Cells.Item(5, 30).Value = 13473.7
$ws.Cells.Item(5, 31).Value = 0
$ws.Cells.Item(5, 32).Value = 13473.7
$ws.Cells.Item(5, 33).Value = 12347.5
$ws.Cells.Item(5, 34).Value = 0.5953595036940154
$ws.Cells.Item(5, 35).Value = 0.8049863183930983
$ws.Cells.Item(5, 36).Value = 0.5741687979539642
$ws.Cells.Item(5, 37).Value = 0.7909182915268134
$ws.Cells.Item(5, 38).Value = 0
$ws.Cells.Item(5, 39).Value = -8.34
$ws.Cells.Item(5, 43).Value = -0.0

# --- Row 6 ---
$ws.Cells.Item(6, 1).Value = "Taiwan"
$ws.Cells.Item(6, 2).Value = "Yuanta Financial Holding Co., Ltd (TSEC:2885)"
$ws.Cells.Item(6, 3).Value = "Financial Svcs. (Non-bank & Insurance)"
$ws.Cells.Item(6, 4).Value = 0.09210000000000002
$ws.Cells.Item(6, 5).Value = 0.119
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 9).Value = 0
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 787.8
$ws.Cells.Item(6, 12).Value = 0.3238243998684643
$ws.Cells.Item(6, 13).Value = 262
$ws.Cells.Item(6, 14).Value = 0.02950948921552064
$ws.Cells.Item(6, 15).Value = 0.3325717187103326
$ws.Cells.Item(6, 16).Value = 262
$ws.Cells.Item(6, 17).Value = 0.02950948921552064
$ws.Cells.Item(6, 18).Value = 0.3325717187103326
$ws.Cells.Item(6, 19).Value = 0
$ws.Cells.Item(6, 20).Value = 0
$ws.Cells.Item(6, 21).Value = 2641.5
$ws.Cells.Item(6, 22).Value = 0.2975164723770907
$ws.Cells.Item(6, 23).Value = 0.1067017011593889
$ws.Cells.Item(6, 24).Value = 0.03042624206767388
$ws.Cells.Item(6, 25).Value = 0.07627545909171499
$ws.Cells.Item(6, 26).Value = 0.1274992269756668
$ws.Cells.Item(6, 27).Value = 0
$ws.Cells.Item(6, 28).Value = 0.02351394358021246
$ws.Cells.Item(6, 29).Value = -0.02351394358021246
$ws.Cells.Item(6, 30).Value = 15814.5
$ws.Cells.Item(6, 31).Value = 0
$ws.Cells.Item(6, 32).Value = 15814.5
$ws.Cells.Item(6, 33).Value = 13173
$ws.Cells.Item(6, 34).Value = 0.6404446604300814
$ws.Cells.Item(6, 35).Value = 0.6338172111962551
$ws.Cells.Item(6, 36).Value = 0.597374328276988
$ws.Cells.Item(6, 37).Value = 0.5904606516448002
$ws.Cells.Item(6, 38).Value = 0
$ws.Cells.Item(6, 39).Value = 0

# --- Row 7 ---
$ws.Cells.Item(7, 1).Value = "Taiwan"
$ws.Cells.Item(7, 2).Value = "IBF Financial Holdings Co., Ltd. (TSEC:2889)"
$ws.Cells.Item(7, 3).Value = "Financial Svcs. (Non-bank & Insurance)"
$ws.Cells.Item(7, 4).Value = 0.112
$ws.Cells.Item(7, 5).Value = 0.151
$ws.Cells.Item(7, 7).Value = 0
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 102.6
$ws.Cells.Item(7, 12).Value = 0.38
$ws.Cells.Item(7, 13).Value = 67.57
$ws.Cells.Item(7, 14).Value = 0.05100007547739451
$ws.Cells.Item(7, 15).Value = 0.6585769980506823
$ws.Cells.Item(7, 16).Value = 63.9
$ws.Cells.Item(7, 17).Value = 0.048230055098498
$ws.Cells.Item(7, 18).Value = 0.6228070175438597
$ws.Cells.Item(7, 19).Value = 3.669999999999995
$ws.Cells.Item(7, 20).Value = 0.05431404469439093
$ws.Cells.Item(7, 21).Value = 111.9
$ws.Cells.Item(7, 22).Value = 0.08445920446826176
$ws.Cells.Item(7, 23).Value = 0.09015025041736227
$ws.Cells.Item(7, 24).Value = 0.06122600878392187
$ws.Cells.Item(7, 25).Value = 0.0289242416334404
$ws.Cells.Item(7, 26).Value = 0.03193841327597036
$ws.Cells.Item(7, 27).Value = 0
$ws.Cells.Item(7, 28).Value = 0.0254097282158323
$ws.Cells.Item(7, 29).Value = -0.0254097282158323
$ws.Cells.Item(7, 30).Value = 8214.9
$ws.Cells.Item(7, 31).Value = 0
$ws.Cells.Item(7, 32).Value = 8214.9
$ws.Cells.Item(7, 33).Value = 8103
$ws.Cells.Item(7, 34).Value = 0.8611186817333697
$ws.Cells.Item(7, 35).Value = 0.8519559445781133
$ws.Cells.Item(7, 36).Value = 0.8594702956119603
$ws.Cells.Item(7, 37).Value = 0.8502177220502597
$ws.Cells.Item(7, 38).Value = 0
$ws.Cells.Item(7, 39).Value = 0

# --- Row 8 ---
$ws.Cells.Item(8, 1).Value = "Taiwan"
$ws.Cells.Item(8, 2).Value = "Hotai Finance Co., Ltd. (TSEC:6592)"
$ws.Cells.Item(8, 3).Value = "Financial Svcs. (Non-bank & Insurance)"
$ws.Cells.Item(8, 7).Value = 0
$ws.Cells.Item(8, 8).Value = 0
$ws.Cells.Item(8, 9).Value = 0
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 11).Value = 82.8
$ws.Cells.Item(8, 12).Value = 0.2141748577340921
$ws.Cells.Item(8, 13).Value = 55.1
$ws.Cells.Item(8, 14).Value = 0.03729020032485111
$ws.Cells.Item(8, 15).Value = 0.6654589371980677
$ws.Cells.Item(8, 16).Value = 55.1
$ws.Cells.Item(8, 17).Value = 0.03729020032485111
$ws.Cells.Item(8, 18).Value = 0.6654589371980677
$ws.Cells.Item(8, 19).Value = 0
$ws.Cells.Item(8, 20).Value = 0
$ws.Cells.Item(8, 21).Value = 18.4
$ws.Cells.Item(8, 22).Value = 0.01245262587980509
$ws.Cells.Item(8, 23).Value = 0.238204833141542
$ws.Cells.Item(8, 24).Value = 0.03780796246838241
$ws.Cells.Item(8, 25).Value = 0.2003968706731596
$ws.Cells.Item(8, 26).Value = 0.1013527684563758
$ws.Cells.Item(8, 27).Value = 0
$ws.Cells.Item(8, 28).Value = 0.02827096076124982
$ws.Cells.Item(8, 29).Value = -0.02827096076124982
$ws.Cells.Item(8, 30).Value = 4196.9
$ws.Cells.Item(8, 31).Value = 0
$ws.Cells.Item(8, 32).Value = 4196.9
$ws.Cells.Item(8, 33).Value = 4178.5
$ws.Cells.Item(8, 34).Value = 0.7396070138338179
$ws.Cells.Item(8, 35).Value = 0.8421253285711419
$ws.Cells.Item(8, 36).Value = 0.7387599229150827
$ws.Cells.Item(8, 37).Value = 0.8415402896098926
$ws.Cells.Item(8, 38).Value = 0.146
$ws.Cells.Item(8, 39).Value = 0.146
$ws.Cells.Item(8, 41).Value = 0
$ws.Cells.Item(8, 43).Value = 0

# Clear cells that existed before but are removed in the updated data
$ws.Cells.Item(5, 40).ClearContents()
$ws.Cells.Item(5, 42).ClearContents()
$ws.Cells.Item(6, 40).ClearContents()
$ws.Cells.Item(6, 42).ClearContents()
# Row 8 inherited AN/AP from the old row 7 (Hotai) via the row-insert shift; the
# updated data no longer has debt_ebitda/net_debt_ebitda for this company.
$ws.Cells.Item(8, 40).ClearContents()
$ws.Cells.Item(8, 42).ClearContents()
